$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates derived from the source diff: refreshed
# DATA_EXTRACCIO timestamps plus a handful of updated meteorological readings.
# Percentage-looking text (e.g. "27%") is forced to stay text by setting the
# NumberFormat to "@" first, otherwise Excel auto-converts it to a numeric
# percentage value, which would not match the original inline-string cell.
$ws.Range("E2").Value = "2026-02-22 19:18:41"
$ws.Range("O2").Value = "6.3 °C"
$ws.Range("E3").Value = "2026-02-22 19:18:44"
$ws.Range("E4").Value = "2026-02-22 19:18:46"
$ws.Range("O4").Value = "12.6 °C"
$ws.Range("E5").Value = "2026-02-22 19:18:49"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "27%"
$ws.Range("O5").Value = "6.3 °C"
$ws.Range("E6").Value = "2026-02-22 19:18:52"
$ws.Range("E7").Value = "2026-02-22 19:18:54"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "58%"
$ws.Range("J7").Value = "1027.4 hPa"
$ws.Range("E8").Value = "2026-02-22 19:18:57"
$ws.Range("E9").Value = "2026-02-22 19:19:00"
$ws.Range("E10").Value = "2026-02-22 19:19:02"
$ws.Range("E11").Value = "2026-02-22 19:19:05"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "64%"
$ws.Range("E12").Value = "2026-02-22 19:19:07"
$ws.Range("O12").Value = "10.0 °C"
$ws.Range("E13").Value = "2026-02-22 19:19:10"
$ws.Range("J13").Value = "1030.3 hPa"
$ws.Range("O13").Value = "6.5 °C"
$ws.Range("E14").Value = "2026-02-22 19:19:12"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "71%"
$ws.Range("O14").Value = "12.2 °C"
$ws.Range("E15").Value = "2026-02-22 19:19:15"
$ws.Range("E16").Value = "2026-02-22 19:19:17"
$ws.Range("O16").Value = "5.4 °C"
$ws.Range("E17").Value = "2026-02-22 19:19:20"
$ws.Range("L17").Value = "34.6 km/h - 267º 18:49 TU"
$ws.Range("O17").Value = "10.2 °C"
$ws.Range("E18").Value = "2026-02-22 19:19:22"
$ws.Range("O18").Value = "10.3 °C"
$ws.Range("E19").Value = "2026-02-22 19:19:25"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "45%"
$ws.Range("E20").Value = "2026-02-22 19:19:28"
$ws.Range("E21").Value = "2026-02-22 19:19:31"
$ws.Range("J21").Value = "1029.1 hPa"
$ws.Range("O21").Value = "9.2 °C"
$ws.Range("E22").Value = "2026-02-22 19:19:33"
$ws.Range("L22").Value = "23.8 km/h - 285º 18:55 TU"
$ws.Range("E23").Value = "2026-02-22 19:19:36"
$ws.Range("E24").Value = "2026-02-22 19:19:39"
$ws.Range("J24").Value = "1029.7 hPa"
$ws.Range("E25").Value = "2026-02-22 19:19:41"
$ws.Range("E26").Value = "2026-02-22 19:19:44"
$ws.Range("N26").Value = "7.1 °C 18:58 TU"
$ws.Range("E27").Value = "2026-02-22 19:19:47"
$ws.Range("O27").Value = "6.8 °C"
$ws.Range("E28").Value = "2026-02-22 19:19:49"
$ws.Range("K28").Value = "15.0 MJ/m2"
$ws.Range("E29").Value = "2026-02-22 19:19:52"
$ws.Range("E30").Value = "2026-02-22 19:19:54"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "71%"
$ws.Range("E31").Value = "2026-02-22 19:19:57"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "63%"
$ws.Range("O31").Value = "14.5 °C"
$ws.Range("E32").Value = "2026-02-22 19:19:59"
$ws.Range("E33").Value = "2026-02-22 19:20:02"
$ws.Range("E34").Value = "2026-02-22 19:20:05"
$ws.Range("O34").Value = "4.5 °C"
$ws.Range("E35").Value = "2026-02-22 19:20:08"
$ws.Range("O35").Value = "11.2 °C"
$ws.Range("E36").Value = "2026-02-22 19:20:10"
$ws.Range("J36").Value = "1027.4 hPa"
$ws.Range("E37").Value = "2026-02-22 19:20:13"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "64%"
$ws.Range("E38").Value = "2026-02-22 19:20:16"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "67%"
$ws.Range("E39").Value = "2026-02-22 19:20:18"
$ws.Range("E40").Value = "2026-02-22 19:20:21"
$ws.Range("J40").Value = "1029.1 hPa"
$ws.Range("E41").Value = "2026-02-22 19:20:24"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "75%"
$ws.Range("E42").Value = "2026-02-22 19:20:27"
$ws.Range("E43").Value = "2026-02-22 19:20:29"
$ws.Range("O43").Value = "9.3 °C"
$ws.Range("E44").Value = "2026-02-22 19:20:31"
$ws.Range("E45").Value = "2026-02-22 19:20:34"
$ws.Range("E46").Value = "2026-02-22 19:20:37"
$ws.Range("J46").Value = "1029.7 hPa"
$ws.Range("O46").Value = "9.1 °C"
